# Apply edits described by the commit diff to Daily_Data.xlsx (Sheet 1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header rename: C1 "RG0_Diff" -> "RG_Diff"
$ws.Range("C1").Value = "RG_Diff"

# 2) Column A holds literal text dates (e.g. "2021/10/1"), not real Excel dates.
#    Force the whole data range to Text format first so re-assigning the date-like
#    strings below is stored as literal text (matches original shared-string layout)
#    instead of being auto-converted into a date serial number.
$ws.Range("A2:A32").NumberFormat = "@"

# 3) Rewrite the date labels (October 2021 -> December 2021) and the Water_Level /
#    RG_Diff measurement columns with the new data.
$ws.Range("A2").Value = "2021/12/1"
$ws.Range("B2").Value = 1.606017
$ws.Range("C2").Value = 2.84248
$ws.Range("A3").Value = "2021/12/2"
$ws.Range("B3").Value = 1.29442
$ws.Range("C3").Value = 5.26993
$ws.Range("A4").Value = "2021/12/3"
$ws.Range("B4").Value = 1.701893
$ws.Range("C4").Value = 12.51669
$ws.Range("A5").Value = "2021/12/4"
$ws.Range("B5").Value = 1.253330285714286
$ws.Range("C5").Value = 0
$ws.Range("A6").Value = "2021/12/5"
$ws.Range("B6").Value = 1.198544
$ws.Range("C6").Value = 0.02434
$ws.Range("A7").Value = "2021/12/6"
$ws.Range("B7").Value = 1.198544
$ws.Range("C7").Value = 0.01687
$ws.Range("A8").Value = "2021/12/7"
$ws.Range("B8").Value = 1.174575
$ws.Range("C8").Value = 0.16428
$ws.Range("A9").Value = "2021/12/8"
$ws.Range("B9").Value = 1.126637
$ws.Range("C9").Value = 0
$ws.Range("A10").Value = "2021/12/9"
$ws.Range("B10").Value = 1.126637
$ws.Range("C10").Value = 3.33038
$ws.Range("A11").Value = "2021/12/10"
$ws.Range("B11").Value = 1.510141
$ws.Range("C11").Value = 2.90854
$ws.Range("A12").Value = "2021/12/11"
$ws.Range("B12").Value = 1.869676
$ws.Range("C12").Value = 2.39
$ws.Range("A13").Value = "2021/12/12"
$ws.Range("B13").Value = 1.582048
$ws.Range("C13").Value = 5.99
$ws.Range("A14").Value = "2021/12/13"
$ws.Range("B14").Value = 1.917614
$ws.Range("C14").Value = 3.21371
$ws.Range("A15").Value = "2021/12/14"
$ws.Range("B15").Value = 1.7738
$ws.Range("C15").Value = 0.7083699999999999
$ws.Range("A16").Value = "2021/12/15"
$ws.Range("B16").Value = 1.246482
$ws.Range("C16").Value = 10.67993
$ws.Range("A17").Value = "2021/12/16"
$ws.Range("B17").Value = 2.041090666666666
$ws.Range("C17").Value = 20.68171
$ws.Range("A18").Value = "2021/12/17"
$ws.Range("B18").Value = 2.396994
$ws.Range("C18").Value = 1.17
$ws.Range("A19").Value = "2021/12/18"
$ws.Range("B19").Value = 1.845707
$ws.Range("C19").Value = 0
$ws.Range("A20").Value = "2021/12/19"
$ws.Range("B20").Value = 1.582048
$ws.Range("C20").Value = 0
$ws.Range("A21").Value = "2021/12/20"
$ws.Range("B21").Value = 1.414265
$ws.Range("C21").Value = 1.39796
$ws.Range("A22").Value = "2021/12/21"
$ws.Range("B22").Value = 1.318389
$ws.Range("C22").Value = 1.05275
$ws.Range("A23").Value = "2021/12/22"
$ws.Range("B23").Value = 1.246482
$ws.Range("C23").Value = 0.25737
$ws.Range("A24").Value = "2021/12/23"
$ws.Range("B24").Value = 1.198544
$ws.Range("C24").Value = 0.03491
$ws.Range("A25").Value = "2021/12/24"
$ws.Range("B25").Value = 1.174575
$ws.Range("C25").Value = 0.57
$ws.Range("A26").Value = "2021/12/25"
$ws.Range("B26").Value = 1.342358
$ws.Range("C26").Value = 2.94
$ws.Range("A27").Value = "2021/12/26"
$ws.Range("B27").Value = 1.318389
$ws.Range("C27").Value = 0.62
$ws.Range("A28").Value = "2021/12/27"
$ws.Range("B28").Value = 1.342358
$ws.Range("C28").Value = 0.47368
$ws.Range("A29").Value = "2021/12/28"
$ws.Range("B29").Value = 1.270451
$ws.Range("C29").Value = 2.13196
$ws.Range("A30").Value = "2021/12/29"
$ws.Range("B30").Value = 1.078699
$ws.Range("C30").Value = 3.58021
$ws.Range("A31").Value = "2021/12/30"
$ws.Range("B31").Value = 1.078699
$ws.Range("C31").Value = 0.0492
$ws.Range("A32").Value = "2021/12/31"
$ws.Range("B32").Value = 1.174575
$ws.Range("C32").Value = 2.50025
